$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 115, shifting existing rows 115-188 down to 119-192.
$ws.Rows.Item(115).Resize(4).Insert()

# Constant column values shared by every data row in this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$prodId    = 100102
$producto  = "Cítricos"
$catId     = 100102005
$categoria = "Naranja"
$unidad    = "$/bandeja 15 kilos granel"
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

# New rows 115-118: data for 2021-09-09 (serial 44438).
$newRows = @(
    @{ Row = 115; Fecha = 44438; Variedad = "Fukumoto";   Calidad = "Primera"; Volumen = 180; Min = 6000; Max = 6500; Prom = 6250; PrecioKg = 417 },
    @{ Row = 116; Fecha = 44438; Variedad = "Fukumoto";   Calidad = "Segunda"; Volumen = 120; Min = 5000; Max = 5500; Prom = 5250; PrecioKg = 350 },
    @{ Row = 117; Fecha = 44438; Variedad = "Navel Late"; Calidad = "Primera"; Volumen = 180; Min = 6000; Max = 6500; Prom = 6250; PrecioKg = 417 },
    @{ Row = 118; Fecha = 44438; Variedad = "Navel Late"; Calidad = "Segunda"; Volumen = 120; Min = 5000; Max = 5500; Prom = 5250; PrecioKg = 350 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
